# Fix header labels on the existing sheets and add a new "PO Forecast"
# sheet containing the forecast series (ds / PO_Forecast / yhat_lower /
# yhat_upper).

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsMonthly = $wb.Worksheets.Item(2)   # "Monthly Trend"

# --- 1) Rename the "Requested quantity" header on both existing sheets ---
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "PO Forecast"

# Header row
$ws.Range("A1").Value = "ds"
$ws.Range("B1").Value = "PO_Forecast"
$ws.Range("C1").Value = "yhat_lower"
$ws.Range("D1").Value = "yhat_upper"

# Reuse the same header formatting (bold, centered, bordered) already used
# on the other two sheets' header rows, instead of building a new style.
$wsWeekly.Range("A1:B1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("C1:D1").PasteSpecial(-4122)   # xlPasteFormats

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$data = @(
    @(45340.99999999999, 37, 6.888807972455351, 67.5259091737632),
    @(45347.99999999999, 37, 7.557423072357192, 66.34502698099578),
    @(45354.99999999999, 36, 7.591429744678538, 66.00400198459923),
    @(45361.99999999999, 35, 4.392126236193677, 62.92930298465507),
    @(45382.99999999999, 33, 2.77912135492751, 61.75942660694165),
    @(45501.99999999999, 21, -7.417719109613281, 49.78538810434504),
    @(45515.99999999999, 19, -10.57605347565959, 49.9541117493597),
    @(45543.99999999999, 16, -14.95063684294491, 44.36755771863358),
    @(45564.99999999999, 14, -15.0179920796298, 45.51187611579082),
    @(45571.99999999999, 13, -15.78322781080995, 42.36443530435493),
    @(45578.99999999999, 13, -17.79247174244934, 40.74390385442265),
    @(45585.99999999999, 12, -16.22354320096573, 39.98599585115335),
    @(45592.99999999999, 11, -18.61980458940881, 37.60871874039444),
    @(45599.99999999999, 10, -19.04738933611745, 39.67250335261501),
    @(45606.99999999999, 10, -20.76041824612881, 37.27286714105362),
    @(45613.99999999999, 9, -20.21777905163977, 37.20773568055503),
    @(45620.99999999999, 8, -20.80484675195376, 39.16907900511303)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# Column A holds dates; copy the existing date-cell format (from the
# "Weekly Quantity" sheet) onto the new column A data range so it reuses
# the same number format style instead of a fresh one.
$wsWeekly.Range("A2").Copy()
$ws.Range("A2:A18").PasteSpecial(-4122)   # xlPasteFormats

$excel.CutCopyMode = 0

# Keep the original active sheet/selection (sheet 1) as it was before the
# new sheet was appended.
$wsWeekly.Activate() | Out-Null
$wsWeekly.Range("A1").Select() | Out-Null
